$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Strip the leading glyph from the "Is Active" status cells
#    (F2 = JOHN DOE, F3 = JANE SMITH, F5 = ALICE BROWN -> "Active";
#     F4 = BOB JOHNSON -> "Inactive")
$ws.Range("F2").Value = "Active"
$ws.Range("F3").Value = "Active"
$ws.Range("F4").Value = "Inactive"
$ws.Range("F5").Value = "Active"

# 2. Data rows should wrap text and left-align (in addition to the
#    existing top alignment) so the now-shorter status text still reads
#    the same as the rest of the body copy.
$body = $ws.Range("A2:G5")
$body.WrapText = $true
$body.HorizontalAlignment = -4131  # xlLeft
$body.VerticalAlignment = -4160    # xlTop

# 3. Columns get one extra character of breathing room.
$ws.Columns("A").ColumnWidth = 14.864285714285714
$ws.Columns("B").ColumnWidth = 14.864285714285714
$ws.Columns("C").ColumnWidth = 18.864285714285714
$ws.Columns("D").ColumnWidth = 11.864285714285714
$ws.Columns("E").ColumnWidth = 15.864285714285716
$ws.Columns("F").ColumnWidth = 10.864285714285714
$ws.Columns("G").ColumnWidth = 9.864285714285714

# 4. Header + data rows grow from 13pt to 14pt.
$ws.Rows("1:5").RowHeight = 14
